$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "id" column was inserted before column A, shifting PLATE/type/color/...
# headers and the existing car row one column to the right (B..H instead of
# A..G). This mirrors the openpyxl/pandas bug described in the commit
# message: an id column got prepended to the sheet.
$ws.Columns.Item(1).Insert()

# Give the new A2 id cell the same look as the header row (bold/centered/
# bordered) by copying the header's format onto it, then write the id value.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 2

# The "second car" record overwrites the plate value that shifted into B2
# (the previous run's "1234567890" plate), recording a brand-new car.
$ws.Range("B2").Value = "second car"

# entry_time / ticket_id reflect this second car's newly generated values.
$ws.Range("G2").Value = 1768578334.572511
$ws.Range("H2").Value = "second car-2602"
